$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

# Copy the style from an existing header cell (F1) so G1/H1 match formatting
$ws.Range("F1").Copy() | Out-Null
$ws.Range("G1:H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Data rows 2-5
for ($r = 2; $r -le 5; $r++) {
    $ws.Cells.Item($r, 7).Value = 0.4788041146331428   # column G
    $ws.Cells.Item($r, 8).Value = 0.997                 # column H
}
